# Delete the "GARAMYCIN 0.1% OINT. 15 GM" line-item row (row 17) from the
# sheet. This shifts all subsequent rows up by one, which Excel handles
# automatically (including adjusting merged cells, the running total in
# column P, and the trailing summary rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(17).Delete()

# Renumber the "م" (item #) column so the sequence stays contiguous
# (1..33) after the row was removed. Data rows now run from row 7 to
# row 39 (previously 7 to 40).
for ($r = 7; $r -le 39; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 6
}

# Update the printed grand total (was 1402.17, minus the removed item's
# price of 22.00). After the row deletion above, this cell is now P40
# (was P41).
$ws.Range("P40").Value = 1380.1700000000001

# Refresh the printed timestamp string that appears near the bottom of the
# report to reflect the new save time (one minute later than before). After
# the row deletion above, this cell is now A41 (was A42).
$ws.Range("A41").Value = "Tuesday, 29 July, 2025 12:53 PM"
